$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link swaps (plain text columns B & C) ---
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'

# --- Price / Volume columns (D & E) ---
# These values must stay stored as TEXT (matching the source feed's
# inline-string cells), but Excel auto-coerces plain numeric-looking
# strings (e.g. "315.74") into real numbers on assignment. Forcing the
# cell to a Text number format first prevents that coercion; clearing
# the format afterwards restores the original (unstyled) cell so only
# the value itself changes.
$priceCells = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","E14","D15","E15","D16","E16","D17","E17","E18","D19","E19","D20","E20","E21","D22","E22","D23","E23","D24","E24","E25","D26","E26","D27","E27","D28","E28","D29","E29","D30","E30","D31","E31","D32","E32","D33","E33","D34","D35","E35","D36","E36","D37","E37","D38","E38","D39","E39","D40","E40","E41","D42","E42","D43","D44","E44","D45","E45","E46","D47","E47","D48","E48","D49","E49","D50","E50","E51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.551.10'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.825.17'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '315.74'
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.5104'
$ws.Range("E7").Value = '  -5.33%  '
$ws.Range("D8").Value = '0.3959'
$ws.Range("E8").Value = '  -1.22%  '
$ws.Range("D9").Value = '0.08252'
$ws.Range("E9").Value = '  +6.18%  '
$ws.Range("D10").Value = '1.115'
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("D11").Value = '41.78'
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("D12").Value = '6.369'
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").Value = '21.21'
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("D15").Value = '7.554'
$ws.Range("E15").Value = '  -1.19%  '
$ws.Range("D16").Value = '1.816.39'
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("D17").Value = '0.00001128'
$ws.Range("E17").Value = '  +3.17%  '
$ws.Range("E18").Value = '  +3.27%  '
$ws.Range("D19").Value = '0.06643'
$ws.Range("E19").Value = '  +0.75%  '
$ws.Range("D20").Value = '17.84'
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").Value = '6.095'
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("D23").Value = '28.585.01'
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").Value = '11.45'
$ws.Range("E24").Value = '  +2.26%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").Value = '21.37'
$ws.Range("E26").Value = '  +2.42%  '
$ws.Range("D27").Value = '156.87'
$ws.Range("E27").Value = '  -0.86%  '
$ws.Range("D28").Value = '2.029.18'
$ws.Range("E28").Value = '  -0.46%  '
$ws.Range("D29").Value = '2.415'
$ws.Range("E29").Value = '  -1.96%  '
$ws.Range("D30").Value = '127.05'
$ws.Range("E30").Value = '  +2.10%  '
$ws.Range("D31").Value = '1.115'
$ws.Range("E31").Value = '  -1.77%  '
$ws.Range("D32").Value = '0.1085'
$ws.Range("E32").Value = '  -3.01%  '
$ws.Range("D33").Value = '5.788'
$ws.Range("E33").Value = '  +1.54%  '
$ws.Range("D34").Value = '3.659'
$ws.Range("D35").Value = '0.07049'
$ws.Range("E35").Value = '  -6.07%  '
$ws.Range("D36").Value = '0.2234'
$ws.Range("E36").Value = '  -0.92%  '
$ws.Range("D37").Value = '0.02358'
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").Value = '5.271'
$ws.Range("E38").Value = '  +1.00%  '
$ws.Range("D39").Value = '8.780'
$ws.Range("E39").Value = '  -2.23%  '
$ws.Range("D40").Value = '0.6365'
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("E41").Value = '  -0.91%  '
$ws.Range("D42").Value = '1.180'
$ws.Range("E42").Value = '  -0.90%  '
$ws.Range("D43").Value = '1.402'
$ws.Range("D44").Value = '13.57'
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").Value = '0.5968'
$ws.Range("E45").Value = '  +0.84%  '
$ws.Range("E46").Value = '  +0.76%  '
$ws.Range("D47").Value = '125.41'
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").Value = '1.996'
$ws.Range("E48").Value = '  -0.42%  '
$ws.Range("D49").Value = '1.194'
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("D50").Value = '0.06937'
$ws.Range("E50").Value = '  +0.32%  '
$ws.Range("E51").Value = '  +4.15%  '

foreach ($addr in $priceCells) {
    $ws.Range($addr).ClearFormats()
}
